$d = $word.ActiveDocument
$bullet = [char]0x2022

# Resolve the 1-based Paragraphs index of the paragraph that contains
# $searchText by using Find (robust to exact paragraph-count assumptions),
# then counting how many paragraph marks precede the hit.
function Get-ParagraphIndexByText($doc, $searchText) {
    $r = $doc.Content
    $ok = $r.Find.Execute($searchText, $false, $false, $false, $false, $false, `
                           $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Text not found: $searchText"
    }
    $null = $r.Expand(4)   # wdParagraph -> grow the hit to its full paragraph
    $pre = $doc.Range(0, $r.Start)
    return $pre.Paragraphs.Count + 1
}

# --- Change 1: the CORE COMPETENCIES section used to list three long detail
#     paragraphs (Research and Analytics / Programming and Development /
#     Data Infrastructure). Collapse them into a single short summary line. ---
$idx1 = Get-ParagraphIndexByText $d "Research and Analytics: Survey Methodology"
$d.Paragraphs.Item($idx1).Range.Text = `
    "Research and Analytics $bullet Programming and Development $bullet Data Infrastructure"
$d.Paragraphs.Item($idx1 + 2).Range.Delete()
$d.Paragraphs.Item($idx1 + 1).Range.Delete()

# --- Change 2: add a new "TECHNICAL SKILLS" section (a Heading2 title plus
#     three detail lines) right after "Built comprehensive survey operations
#     platform from RFP through deployment" and before the closing
#     "For a more detailed..." paragraph. ---
$idxBuilt = Get-ParagraphIndexByText $d "Built comprehensive survey operations platform"
$anchor = $d.Paragraphs.Item($idxBuilt)
$r = $anchor.Range
$r.Collapse(0)
$r.InsertAfter( `
    "TECHNICAL SKILLS`r" + `
    "RESEARCH AND ANALYTICS Survey Methodology; Statistical Analysis; Geospatial Analysis; Data Visualization`r" + `
    "PROGRAMMING AND DEVELOPMENT Python; JVM Languages; Web Technologies; Database Languages`r" + `
    "DATA INFRASTRUCTURE Cloud Platforms; Big Data; Databases; Geospatial`r")

$d.Paragraphs.Item($idxBuilt + 1).Style = "Heading 2"
